$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Remove the data row for "DI MAURO SALVATORE" (row 7). Excel shifts every
# row below it up by one, and the now-unreferenced shared string is dropped
# automatically on save.
$ws.Rows.Item(7).Delete()

# The autofilter range shrank along with the data (rows 2:25 -> 2:24), so
# rebuild it over the new extent A1:D24.
$ws.AutoFilterMode = $false
$ws.Range("A1:D24").AutoFilter()

# Keep the workbook-level _FilterDatabase defined name in sync with the
# refreshed autofilter range.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Foglio1!`$A`$1:`$D`$24"
    }
}

# Restore the default 100% zoom (was 110%) and move the active selection
# to A13, matching the saved view state.
$excel.ActiveWindow.Zoom = 100
$ws.Range("A13").Select()
